$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(9, 8).Value = 1127
$ws.Cells.Item(9, 9).Value = 1291
$ws.Cells.Item(9, 10).Value = 471
$ws.Cells.Item(9, 11).Value = 1291
$ws.Cells.Item(9, 12).Value = 471
$ws.Cells.Item(9, 13).Value = -1122
$ws.Cells.Item(9, 14).Value = -809
$ws.Cells.Item(98, 8).Value = 2661.3215
$ws.Cells.Item(98, 9).Value = 2522.2083
$ws.Cells.Item(98, 11).Value = 2522.2083
$ws.Cells.Item(98, 13).Value = -1024.2083
$ws.Cells.Item(105, 8).Value = 25000
$ws.Cells.Item(105, 10).Value = 25000
$ws.Cells.Item(105, 12).Value = 25000
$ws.Cells.Item(105, 14).Value = -31988
$ws.Cells.Item(113, 8).Value = 12842.857
$ws.Cells.Item(113, 9).Value = 16780
$ws.Cells.Item(113, 11).Value = 16780
$ws.Cells.Item(113, 13).Value = -13526
$ws.Cells.Item(122, 8).Value = 2661.3215
$ws.Cells.Item(122, 9).Value = 2522.2083
$ws.Cells.Item(122, 11).Value = 7566.624899999999
$ws.Cells.Item(122, 13).Value = -5116.624899999999
$ws.Cells.Item(135, 8).Value = 1502.7894
$ws.Cells.Item(135, 9).Value = 1597.25
$ws.Cells.Item(135, 11).Value = 14375.25
$ws.Cells.Item(135, 13).Value = -11840.25
$ws.Cells.Item(137, 8).Value = 11601.55
$ws.Cells.Item(137, 9).Value = 1496.8
$ws.Cells.Item(137, 10).Value = 21706.3
$ws.Cells.Item(137, 11).Value = 4490.4
$ws.Cells.Item(137, 12).Value = 65118.89999999999
$ws.Cells.Item(137, 13).Value = -1940.4
$ws.Cells.Item(137, 14).Value = -70218.89999999999
$ws = $wb.Worksheets.Item("ARM")
$ws.Cells.Item(74, 8).Value = 15051.324
$ws.Cells.Item(74, 9).Value = 1798.909
$ws.Cells.Item(74, 10).Value = 34488.2
$ws.Cells.Item(74, 11).Value = 1798.909
$ws.Cells.Item(74, 12).Value = 34488.2
$ws.Cells.Item(74, 13).Value = -924.9090000000001
$ws.Cells.Item(74, 14).Value = -36236.2
$ws.Cells.Item(77, 8).Value = 15051.324
$ws.Cells.Item(77, 9).Value = 1798.909
$ws.Cells.Item(77, 10).Value = 34488.2
$ws.Cells.Item(77, 11).Value = 8994.545
$ws.Cells.Item(77, 12).Value = 172441
$ws.Cells.Item(77, 13).Value = -4626.545
$ws.Cells.Item(77, 14).Value = -181177
$ws.Cells.Item(110, 8).Value = 6878.4443
$ws.Cells.Item(110, 9).Value = 7613.25
$ws.Cells.Item(110, 10).Value = 1000
$ws.Cells.Item(110, 11).Value = 7613.25
$ws.Cells.Item(110, 12).Value = 1000
$ws.Cells.Item(110, 13).Value = -5568.25
$ws.Cells.Item(110, 14).Value = -5090
$ws.Cells.Item(139, 8).Value = 90614.92
$ws.Cells.Item(139, 10).Value = 90614.92
$ws.Cells.Item(139, 12).Value = 90614.92
$ws.Cells.Item(139, 14).Value = -100894.92
$ws = $wb.Worksheets.Item("BSM")
$ws.Cells.Item(99, 8).Value = 8738.571
$ws.Cells.Item(99, 9).Value = 9279.031999999999
$ws.Cells.Item(99, 10).Value = 4550
$ws.Cells.Item(99, 11).Value = 9279.031999999999
$ws.Cells.Item(99, 12).Value = 4550
$ws.Cells.Item(99, 13).Value = -7781.031999999999
$ws.Cells.Item(99, 14).Value = -7546
$ws.Cells.Item(105, 8).Value = 1812
$ws.Cells.Item(105, 9).Value = 1874.4
$ws.Cells.Item(105, 10).Value = 1500
$ws.Cells.Item(105, 11).Value = 1874.4
$ws.Cells.Item(105, 12).Value = 1500
$ws.Cells.Item(105, 13).Value = -127.4000000000001
$ws.Cells.Item(105, 14).Value = -4994
$ws.Cells.Item(107, 8).Value = 1659.0682
$ws.Cells.Item(107, 9).Value = 1703
$ws.Cells.Item(107, 11).Value = 1703
$ws.Cells.Item(107, 13).Value = 217
$ws.Cells.Item(134, 8).Value = 52841.6
$ws.Cells.Item(134, 9).Value = 68528.336
$ws.Cells.Item(134, 10).Value = 29311.5
$ws.Cells.Item(134, 11).Value = 205585.008
$ws.Cells.Item(134, 12).Value = 87934.5
$ws.Cells.Item(134, 13).Value = -203050.008
$ws.Cells.Item(134, 14).Value = -93004.5
$ws = $wb.Worksheets.Item("CRP")
$ws.Cells.Item(31, 8).Value = 11491
$ws.Cells.Item(31, 9).Value = 947.5238000000001
$ws.Cells.Item(31, 11).Value = 947.5238000000001
$ws.Cells.Item(31, 13).Value = -652.5238000000001
$ws.Cells.Item(34, 8).Value = 11491
$ws.Cells.Item(34, 9).Value = 947.5238000000001
$ws.Cells.Item(34, 11).Value = 947.5238000000001
$ws.Cells.Item(34, 13).Value = -745.5238000000001
$ws.Cells.Item(107, 8).Value = 1021.2917
$ws.Cells.Item(107, 9).Value = 962.6842
$ws.Cells.Item(107, 11).Value = 962.6842
$ws.Cells.Item(107, 13).Value = 957.3158
$ws = $wb.Worksheets.Item("CUL")
$ws.Cells.Item(92, 8).Value = 302.91666
$ws.Cells.Item(92, 9).Value = 316.57144
$ws.Cells.Item(92, 10).Value = 283.8
$ws.Cells.Item(92, 11).Value = 949.71432
$ws.Cells.Item(92, 12).Value = 851.4000000000001
$ws.Cells.Item(92, 13).Value = 298.28568
$ws.Cells.Item(92, 14).Value = -3347.4
$ws.Cells.Item(107, 8).Value = 806.9167
$ws.Cells.Item(107, 9).Value = 663.3333
$ws.Cells.Item(107, 10).Value = 950.5
$ws.Cells.Item(107, 11).Value = 1989.9999
$ws.Cells.Item(107, 12).Value = 2851.5
$ws.Cells.Item(107, 13).Value = -69.99990000000003
$ws.Cells.Item(107, 14).Value = -6691.5
$ws = $wb.Worksheets.Item("GSM")
$ws.Cells.Item(68, 8).Value = 329727.34
$ws.Cells.Item(68, 9).Value = 49999
$ws.Cells.Item(68, 10).Value = 469591.5
$ws.Cells.Item(68, 11).Value = 49999
$ws.Cells.Item(68, 12).Value = 469591.5
$ws.Cells.Item(68, 13).Value = -49188
$ws.Cells.Item(68, 14).Value = -471213.5
$ws.Cells.Item(71, 8).Value = 329727.34
$ws.Cells.Item(71, 9).Value = 49999
$ws.Cells.Item(71, 10).Value = 469591.5
$ws.Cells.Item(71, 11).Value = 149997
$ws.Cells.Item(71, 12).Value = 1408774.5
$ws.Cells.Item(71, 13).Value = -145941
$ws.Cells.Item(71, 14).Value = -1416886.5
$ws.Cells.Item(122, 8).Value = 4308.9546
$ws.Cells.Item(122, 9).Value = 4564.2144
$ws.Cells.Item(122, 10).Value = 3862.25
$ws.Cells.Item(122, 11).Value = 13692.6432
$ws.Cells.Item(122, 12).Value = 11586.75
$ws.Cells.Item(122, 13).Value = -11242.6432
$ws.Cells.Item(122, 14).Value = -16486.75
$ws.Cells.Item(123, 8).Value = 58674
$ws.Cells.Item(123, 10).Value = 58674
$ws.Cells.Item(123, 12).Value = 58674
$ws.Cells.Item(123, 14).Value = -63574
$ws = $wb.Worksheets.Item("LTW")
$ws.Cells.Item(55, 8).Value = 2216.8235
$ws.Cells.Item(55, 9).Value = 1832.7778
$ws.Cells.Item(55, 11).Value = 1832.7778
$ws.Cells.Item(55, 13).Value = -1659.7778
$ws.Cells.Item(61, 8).Value = 3292.1428
$ws.Cells.Item(61, 9).Value = 3007.5
$ws.Cells.Item(61, 11).Value = 3007.5
$ws.Cells.Item(61, 13).Value = -2805.5
$ws.Cells.Item(82, 8).Value = 6351.5
$ws.Cells.Item(85, 8).Value = 6351.5
$ws.Cells.Item(106, 8).Value = 18768.428
$ws.Cells.Item(106, 10).Value = 18768.428
$ws.Cells.Item(106, 12).Value = 18768.428
$ws.Cells.Item(106, 14).Value = -21292.428
$ws.Cells.Item(113, 8).Value = 3292.1428
$ws.Cells.Item(113, 9).Value = 3007.5
$ws.Cells.Item(113, 11).Value = 3007.5
$ws.Cells.Item(113, 13).Value = -837.5
$ws.Cells.Item(133, 8).Value = 69999
$ws.Cells.Item(133, 10).Value = 69999
$ws.Cells.Item(133, 12).Value = 69999
$ws.Cells.Item(133, 14).Value = -75059
$ws = $wb.Worksheets.Item("WVR")
$ws.Cells.Item(62, 8).Value = 29378.166
$ws.Cells.Item(62, 9).Value = 37534.5
$ws.Cells.Item(62, 10).Value = 25300
$ws.Cells.Item(62, 11).Value = 37534.5
$ws.Cells.Item(62, 12).Value = 25300
$ws.Cells.Item(62, 13).Value = -36910.5
$ws.Cells.Item(62, 14).Value = -26548
$ws.Cells.Item(65, 8).Value = 29378.166
$ws.Cells.Item(65, 9).Value = 37534.5
$ws.Cells.Item(65, 10).Value = 25300
$ws.Cells.Item(65, 11).Value = 187672.5
$ws.Cells.Item(65, 12).Value = 126500
$ws.Cells.Item(65, 13).Value = -184552.5
$ws.Cells.Item(65, 14).Value = -132740
$ws.Cells.Item(107, 8).Value = 1504.6364
$ws.Cells.Item(107, 9).Value = 1761.2222
$ws.Cells.Item(107, 10).Value = 350
$ws.Cells.Item(107, 11).Value = 5283.6666
$ws.Cells.Item(107, 12).Value = 1050
$ws.Cells.Item(107, 13).Value = -3363.6666
$ws.Cells.Item(107, 14).Value = -4890
$ws.Cells.Item(126, 8).Value = 8610.777
$ws.Cells.Item(126, 9).Value = 8610.777
$ws.Cells.Item(126, 11).Value = 25832.331
$ws.Cells.Item(126, 13).Value = -23362.331
$ws.Cells.Item(136, 8).Value = 301761.03
$ws.Cells.Item(136, 10).Value = 1051458.4
$ws.Cells.Item(136, 12).Value = 3154375.2
$ws.Cells.Item(136, 14).Value = -3159475.2
